$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '91.725.66'
Set-TextCell 'E2' '  +0.61%  '
Set-TextCell 'D3' '3.128.55'
Set-TextCell 'E3' '  +1.19%  '
Set-TextCell 'E4' '  +0.06%  '
Set-TextCell 'D5' '246.33'
Set-TextCell 'E5' '  +0.35%  '
Set-TextCell 'D6' '617.84'
Set-TextCell 'E6' '  -0.40%  '
Set-TextCell 'E7' '  -3.50%  '
Set-TextCell 'D8' '0.385'
Set-TextCell 'E8' '  +4.40%  '
Set-TextCell 'E9' '  -0.02%  '
Set-TextCell 'D10' '3.123.84'
Set-TextCell 'E10' '  +1.10%  '
Set-TextCell 'D11' '0.738'
Set-TextCell 'E11' '  -1.76%  '
Set-TextCell 'E12' '  +1.16%  '
Set-TextCell 'E13' '  +0.33%  '
Set-TextCell 'B14' 'Avalanche'
Set-TextCell 'C14' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 'D14' '34.94'
Set-TextCell 'E14' '  -1.80%  '
Set-TextCell 'B15' 'Toncoin'
Set-TextCell 'C15' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D15' '5.61'
Set-TextCell 'E15' '  +2.33%  '
Set-TextCell 'D16' '91.529.73'
Set-TextCell 'E16' '  +0.43%  '
Set-TextCell 'D17' '3.702.10'
Set-TextCell 'E17' '  +0.94%  '
Set-TextCell 'D18' '3.146.31'
Set-TextCell 'E18' '  +1.58%  '
Set-TextCell 'D19' '3.73'
Set-TextCell 'E19' '  +1.43%  '
Set-TextCell 'D20' '14.87'
Set-TextCell 'E20' '  +1.69%  '
Set-TextCell 'D21' '5.81'
Set-TextCell 'E21' '  -0.33%  '
Set-TextCell 'E22' '  +3.73%  '
Set-TextCell 'D23' '447.77'
Set-TextCell 'E23' '  +0.58%  '
Set-TextCell 'E24' '  -4.60%  '
Set-TextCell 'D25' '5.87'
Set-TextCell 'E25' '  +5.20%  '
Set-TextCell 'D26' '88.15'
Set-TextCell 'E26' '  -3.44%  '
Set-TextCell 'D27' '11.78'
Set-TextCell 'E27' '  -1.18%  '
Set-TextCell 'E28' '  +0.66%  '
Set-TextCell 'D29' '0.145'
Set-TextCell 'E29' '  +30.07%  '
Set-TextCell 'E30' '  +0.01%  '
Set-TextCell 'D31' '0.237'
Set-TextCell 'E31' '  -3.66%  '
Set-TextCell 'E32' '  -10.38%  '
Set-TextCell 'E33' '  +3.82%  '
Set-TextCell 'E34' '  +0.37%  '
Set-TextCell 'E35' '  -0.79%  '
Set-TextCell 'D36' '7.79'
Set-TextCell 'E36' '  +0.56%  '
Set-TextCell 'D37' '26.28'
Set-TextCell 'E37' '  -1.02%  '
Set-TextCell 'D38' '4.16'
Set-TextCell 'E38' '  -0.71%  '
Set-TextCell 'E39' '  +1.16%  '
Set-TextCell 'D40' '492.38'
Set-TextCell 'E40' '  -0.69%  '
Set-TextCell 'E41' '  +0.71%  '
Set-TextCell 'E42' '  +4.70%  '
Set-TextCell 'D43' '3.42'
Set-TextCell 'E43' '  -5.75%  '
Set-TextCell 'E44' '  +0.37%  '
Set-TextCell 'E45' '  -0.04%  '
Set-TextCell 'D46' '157.43'
Set-TextCell 'E46' '  +2.31%  '
Set-TextCell 'D47' '0.707'
Set-TextCell 'E47' '  +1.62%  '
Set-TextCell 'E48' '  +0.50%  '
Set-TextCell 'E49' '  +0.98%  '
Set-TextCell 'D50' '44.17'
Set-TextCell 'E50' '  -0.05%  '
Set-TextCell 'D51' '4.40'
Set-TextCell 'E51' '  -2.04%  '
